# Add new crypto-related words to the word list (translated via Google Translate
# per the commit message) and populate the next four rows in column A only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWords = @("ath", "whale", "ballena", "nft")

$startRow = 32
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newWords[$i]
}

# Move the view/selection down to reflect the newly added rows.
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("A38").Select()
